$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sat Feb 17 22:54:03 EST 2024"
$ws.Range("B3").Value = "Sat Feb 17 22:54:17 EST 2024"
$ws.Range("B5").Value = "Sat Feb 17 22:54:29 EST 2024"
$ws.Range("B6").Value = "Sat Feb 17 22:54:42 EST 2024"
$ws.Range("B7").Value = "Sat Feb 17 22:54:55 EST 2024"
